# Apply the "Add strings to relevant spreadsheets" edit.
#
# Summary of the target change (VIE site-ui-strings translation workbook):
#   Row 62 (old "Results for this location…" / "Kết quả cho địa điểm này…")
#   is replaced so that column A now holds the English string
#   "Results for this location" (no ellipsis) while column B keeps the
#   existing Vietnamese translation "Kết quả cho địa điểm này…".
#   Two brand new rows are appended at the bottom with new English/Vietnamese
#   string pairs, and the old "leftover" blank row 65 (which only carried
#   left-over formatting) becomes a real data row.
#   The old bespoke borderless/Times-New-Roman formatting that rows 61-65
#   used is dropped so the rows fall back to the sheet's normal column
#   formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 62: swap the English text for the non-ellipsis variant; Vietnamese
# column keeps its existing translation, just shifted up onto this row.
$ws.Range("A62").Value = "Results for this location"
$ws.Range("B62").Value = "Kết quả cho địa điểm này…"

# New row 65: "What to Expect at This Location" pair.
$ws.Range("A65").Value = "What to Expect at This Location"
$ws.Range("B65").Value = "Điều Gì Có Thể Xảy Ra tại Địa Điểm Này"

# New row 66: "Getting results for your location…" pair.
$ws.Range("A66").Value = "Getting results for your location…"
$ws.Range("B66").Value = "Xem kết quả cho địa điểm của bạn…"

# Drop the old one-off formatting (Times New Roman / no border) on rows
# 61-66 so they pick up the sheet's normal default formatting again.
$ws.Range("A61:B66").Style = "Normal"

# Match the updated selection left behind in the saved worksheet.
$ws.Range("A61:B66").Select()
